# support nan values
# Fill in newly-collected row-2 data, clear the now-missing (NaN) values
# for rows 3-7, unhide/resize column E, and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: new volunteer data that was previously blank ---
$ws.Range("C2").Value = "שחר שטוקהמר"
$ws.Range("D2").Value = 318162112
$ws.Range("E2").Value = "שחר שטוקהמר"
$ws.Range("F2").Value = "כגדהדג"
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = "fgdfgdfg"

# --- Rows 3-7: these records no longer have the extra (NaN) fields ---
$ws.Range("G3:G7").ClearContents()
$ws.Range("J3:P7").ClearContents()

# --- Column E: unhide and resize to fit the new content ---
$ws.Columns("E").Hidden = $false
$ws.Columns("E").ColumnWidth = 19.2

# --- Update the current selection ---
$ws.Range("E2").Select()
